# Edit script: add "Focus Group " before "Script" in the heading,
# and append "Focus Group Notes" section with three paragraphs of notes
# at the end of the document (before the _GoBack bookmark).

$d = $word.ActiveDocument

# --- Hunk 1: split the "Script" heading run into "Focus Group " + "Script" ---
$r1 = $d.Content
$r1.Find.Execute("Script", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $r1.Find.Found) {
    throw "Could not find 'Script' heading to update"
}
$r1.Collapse(1)
$xml1 = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="5C039D51" w14:textId="0F70942F" w:rsidR="00EF7766" w:rsidRDefault="00986FB3" w:rsidP="00986FB3"><w:pPr><w:spacing w:line="480" w:lineRule="auto"/><w:jc w:val="center"/><w:rPr><w:sz w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:lastRenderedPageBreak/><w:t xml:space="preserve">Focus Group </w:t></w:r><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t>Script</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$r1.InsertXML($xml1)

# --- Hunk 2: append the "Focus Group Notes" section after the closing line ---
$r2 = $d.Content
$r2.Find.Execute("Thank you again for taking the time to answer our questions today.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $r2.Find.Found) {
    throw "Could not find closing line to append notes after"
}
$r2.Collapse(1)
$xml2 = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="2E984C5D" w14:textId="5B8A7C21" w:rsidR="00413221" w:rsidRPr="002209FA" w:rsidRDefault="00413221" w:rsidP="00961E3C"><w:pPr><w:spacing w:line="480" w:lineRule="auto"/><w:ind w:firstLine="720"/><w:rPr><w:sz w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t>Thank you again for taking the time to answer our questions today.</w:t></w:r></w:p><w:p><w:pPr><w:jc w:val="center"/><w:rPr><w:sz w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:br w:type="page"/></w:r><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:lastRenderedPageBreak/><w:t>Focus Group Notes</w:t></w:r></w:p><w:p><w:pPr><w:spacing w:line="360" w:lineRule="auto"/><w:ind w:firstLine="720"/><w:rPr><w:sz w:val="24"/></w:rPr></w:pPr><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t>As a result of</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve"> our lack of development progress, there was not a working demonstration we were able to give. We did, however, show interface mockups and give focus group participants an idea of what we hope Xpendit will be in the future.</w:t></w:r></w:p><w:p><w:pPr><w:spacing w:line="360" w:lineRule="auto"/><w:ind w:firstLine="720"/><w:rPr><w:sz w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t>Each of the participants thought that Xpendit would be a useful app conceptually. One participant noted specifically how Xpendit would be useful for tracking which member of an office owes how much money for a shared purchase of snacks.</w:t></w:r></w:p><w:p><w:pPr><w:spacing w:line="360" w:lineRule="auto"/><w:ind w:firstLine="720"/><w:rPr><w:sz w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve">In terms of user interface, each participant enjoyed </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t>Xpendit’s</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve"> focus on having a dark theme. Not all participants liked the light blue and pink color scheme employed by Xpendit currently, and stated that they would enjoy having the ability to change the primary and accent colors from the default colors they have now.</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$r2.InsertXML($xml2)
